# Generate Report for Handback
# Mark the handed-back rows (7e5ed80a-...) as complete instead of pending,
# and stamp the actual handback datetime for each locale.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("C3").Value = "Handed back: in sync with en-US"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Handed back: in sync with en-US"
$wsZhCn.Range("H2").Value = "2016-03-24 11:01:43"
$wsZhCn.Range("H3").Value = "2016-03-24 11:01:43"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Handed back: in sync with en-US"
$wsDeDe.Range("H2").Value = "2016-03-24 11:02:02"
$wsDeDe.Range("H3").Value = "2016-03-24 11:02:02"
